$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.0003714022599530242
$ws.Range("C2").Value = 0.004309184025731883
$ws.Range("D2").Value = 3.082599426703578
$ws.Range("E2").Value = 246.9852506941017
$ws.Range("G2").Value = 250.0725307070909

# Row 3
$ws.Range("B3").Value = 1.505614041169197
$ws.Range("C3").Value = 1.65323645889881
$ws.Range("D3").Value = 0.7127328510149897
$ws.Range("E3").Value = 0.4998867070740569
$ws.Range("G3").Value = 4.371470058157054
